# Word COM-interop script: expand the MyAnalytics "fitness tracker" blurb
# with the new copy, and move the "_GoBack" bookmark from the end of the
# "Join us on ... Microsoft MyAnalytics." paragraph to the end of the
# blurb paragraph we just rewrote.

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it currently sits after
#    "Microsoft MyAnalytics.").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Locate the "fitness tracker for the workplace" paragraph and cut
#    everything from "at work" onward, then rebuild it run-by-run with
#    the new copy so the text matches the updated announcement.
# ------------------------------------------------------------------
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*fitness tracker for the workplace*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find the 'fitness tracker' paragraph"
}

$pRange = $target.Range
$pStart = $pRange.Start
$pEnd = $pRange.End

$fullText = $pRange.Text
$cutOffset = $fullText.IndexOf("at work") + ("at work").Length
$cutPoint = $pStart + $cutOffset

# Drop the old trailing text (", then set goals ... support your goals.")
# but keep the paragraph mark.
$tail = $d.Range($cutPoint, $pEnd - 1)
$tail.Text = ""

# Re-insert the new continuation as a sequence of inserts (mirrors how
# the edit was authored as several distinct runs).
$pos = $cutPoint

$segments = @(
    " ",
    "and ",
    "find ways to be more productive",
    ". ",
    "You" + [char]0x2019 + "ll get insight into how much high-quality focus time you get each week, how much you collaborate with ",
    "coworkers",
    " after-hours, and much more",
    ". You" + [char]0x2019 + "ll ",
    "also get tips ",
    "directly in Outlook to support your goals."
)

foreach ($seg in $segments) {
    $r = $d.Range($pos, $pos)
    $r.InsertAfter($seg)
    $pos = $pos + $seg.Length
}

# ------------------------------------------------------------------
# 3) Add the "_GoBack" bookmark back at the end of this paragraph
#    (collapsed range, right before the paragraph mark).
# ------------------------------------------------------------------
$newPRange = $target.Range
$bookmarkPoint = $newPRange.End - 1
$bmRange = $d.Range($bookmarkPoint, $bookmarkPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
